$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 153"
$ws.Range("D11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 152"
$ws.Range("I11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 151"
